$d = $word.ActiveDocument

# Locate the paragraph that contains the "Requisitos" section's requirement
# line ("LOT2013: Engenharia Bioquímica I (Requisito fraco)"). The three
# paragraphs that immediately follow it (a blank paragraph, the
# "Ver no Jupiter..." line, and the "© 2020 ..." footer line) were removed
# by the site rebuild, so we delete that whole span in one shot.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOT2013: Engenharia Bioqu*mica I (Requisito fraco)*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $afterIndex = $target.Index + 1

    $first = $d.Paragraphs($afterIndex)
    $last = $d.Paragraphs($afterIndex + 2)

    $start = $first.Range.Start
    $end = $last.Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
